# "logi uuendus (uus nadal)" - add a new week sheet ("Nädal 4") in front of the
# existing weeks, pre-filled with the first day's entry, and make it active.

$wb = $excel.ActiveWorkbook

# The leftmost sheet ("Nädal 3") is the template for a week log - duplicate it
# and place the copy before it so the new sheet becomes the first tab.
$template = $wb.Worksheets.Item(1)
$template.Copy($template)

$newWeek = $wb.Worksheets.Item(1)
$newWeek.Name = "Nädal 4"

# Wipe the copied sample data (rows 8-18 plus the first data row 7) so only
# the new week's own entries remain.
$newWeek.Range("B8:J18").ClearContents()
$newWeek.Range("B7:J7").ClearContents()

# Header date for the new week.
$newWeek.Range("G4").Value = 43885

# First logged entry for the new week.
$newWeek.Range("B7").Value = 43879
$newWeek.Range("C7").Value = 0.33333333333333331
$newWeek.Range("D7").Value = 0.39583333333333331
$newWeek.Range("E7").Value = "-"
$newWeek.Range("F7").Value = 90
$newWeek.Range("G7").Value = "Loeng"
$newWeek.Range("I7").Value = "x"

# Make the new sheet active/selected with the same cell selection the author
# left it in.
$newWeek.Activate()
$newWeek.Range("I7").Select()
